# Append a new data row (row 88) to each of the 4 worksheets, mirroring
# the existing row-87 layout/formatting, and extend each sheet's used range.

$wb = $excel.ActiveWorkbook

# Per-sheet values for the new row 88 (columns A..I).
# G (the big "ID_DEC" number) is supplied as a plain digit string so Excel
# parses it into a double with full precision (scientific-notation literals
# aren't supported by the script parser, and naive float math loses
# precision on round numbers like ...000000000).
$rowsData = @{
    1 = @{
        A = 45874.46114583333
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x20"
        E = "0x07"
        F = 400
        G = "568631262647113000000000"
        H = 288
        I = 7
    }
    2 = @{
        A = 45874.46114583333
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x2C"
        E = "0x19"
        F = 380
        G = "568432987514711000000000"
        H = 300
        I = 25
    }
    3 = @{
        A = 45874.46114583333
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x5F"
        E = "0x15"
        F = 110
        G = "568631262647113000000000"
        H = 95
        I = 15
    }
    4 = @{
        A = 45874.46114583333
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x75"
        E = "0x9"
        F = 130
        G = "568631262647113000000000"
        H = 117
        I = 9
    }
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$i]

    $srcRow = 87
    $dstRow = 88

    # Match column A's existing date/time number format (style) before
    # writing its value, so the new row renders the same as the rest.
    $dateFmt = $ws.Cells.Item($srcRow, 1).NumberFormat()
    $ws.Cells.Item($dstRow, 1).NumberFormat = $dateFmt
    $ws.Cells.Item($dstRow, 1).Value = $data.A

    $ws.Cells.Item($dstRow, 2).Value = $data.B
    $ws.Cells.Item($dstRow, 3).Value = $data.C
    $ws.Cells.Item($dstRow, 4).Value = $data.D
    $ws.Cells.Item($dstRow, 5).Value = $data.E
    $ws.Cells.Item($dstRow, 6).Value = $data.F
    $ws.Cells.Item($dstRow, 7).Value = $data.G
    $ws.Cells.Item($dstRow, 8).Value = $data.H
    $ws.Cells.Item($dstRow, 9).Value = $data.I
}
